$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the existing "Hipervinculo" cell style (currently applied to C2/C3)
# so it can be re-applied after Hyperlinks.Add() touches formatting.
$linkStyle = $ws.Range("C2").Style

# Drop the old hyperlink(s) on the sheet; they'll be rebuilt below with the
# new targets for C2 and the brand-new one for C3.
$ws.Cells.Hyperlinks.Delete()

# --- Row 2: Gigantosaurio photo link changes to the new uploaded image ---
$ws.Range("C2").Value = "https://github.com/SergioGerman412/30DayChartChallenge/blob/main/Day19-Dinosaurs/Im%C3%A1genes/gigantosaurio.png?raw=true"

# --- Row 3: new dinosaur entry (Carcharodontosaurus) ---
$ws.Range("C3").Value = "https://github.com/SergioGerman412/30DayChartChallenge/blob/main/Day19-Dinosaurs/Imágenes/carcharodontosaurus.jpg?raw=true"
$ws.Range("A3").Value = "Carcharodontosaurus"
$ws.Range("B3").Value = "Gigantesco terópodo carnívoro que vivió durante el período Cretácico en lo que hoy es África. Se estima que llegaba a medir hasta 13 metros de largo y poseía enormes dientes serrados, similar al Tyrannosaurus rex."

# Rebuild the hyperlinks pointing at the photo cells
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/SergioGerman412/30DayChartChallenge/blob/main/Day19-Dinosaurs/Im%C3%A1genes/gigantosaurio.png?raw=true")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/SergioGerman412/30DayChartChallenge/blob/main/Day19-Dinosaurs/Imágenes/carcharodontosaurus.jpg?raw=true")

# Restore the hyperlink cell style that Hyperlinks.Add() reformatted
$ws.Range("C2").Style = $linkStyle
$ws.Range("C3").Style = $linkStyle

# The saved selection moves back to a single cell
$ws.Range("C9").Select() | Out-Null
